$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4129.4
$ws.Range("I32").Value = 4932.3335
$ws.Range("J32").Value = 2925
$ws.Range("K32").Value = 4932.3335
$ws.Range("L32").Value = 2925
$ws.Range("M32").Value = -4606.3335
$ws.Range("N32").Value = -3577

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H96").Value = 300
$ws.Range("J96").Value = 300
$ws.Range("L96").Value = 900
$ws.Range("N96").Value = -3646

$ws.Range("H100").Value = 3336100
$ws.Range("I100").Value = 3336100
$ws.Range("K100").Value = 3336100
$ws.Range("M100").Value = -3335559

$ws.Range("H103").Value = 2383.3333
$ws.Range("I103").Value = 2350
$ws.Range("J103").Value = 2450
$ws.Range("K103").Value = 7050
$ws.Range("L103").Value = 7350
$ws.Range("M103").Value = -6464
$ws.Range("N103").Value = -8522

$ws.Range("H135").Value = 1214.1428
$ws.Range("J135").Value = 1061.1666
$ws.Range("L135").Value = 9550.499400000001
$ws.Range("N135").Value = -14620.4994

$ws.Range("H138").Value = 5959.93
$ws.Range("J138").Value = 6573
$ws.Range("L138").Value = 19719
$ws.Range("N138").Value = -29999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2429
$ws.Range("I32").Value = 2388.3215
$ws.Range("J32").Value = 2998.5
$ws.Range("K32").Value = 2388.3215
$ws.Range("L32").Value = 2998.5
$ws.Range("M32").Value = -2101.3215
$ws.Range("N32").Value = -3572.5

$ws.Range("H45").Value = 3009.5
$ws.Range("I45").Value = 3009.5
$ws.Range("K45").Value = 3009.5
$ws.Range("M45").Value = -2632.5

$ws.Range("H97").Value = 1999.6
$ws.Range("J97").Value = 1832.6666
$ws.Range("L97").Value = 1832.6666
$ws.Range("N97").Value = -2824.6666

$ws.Range("H102").Value = 1249
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3289.7058
$ws.Range("I20").Value = 3614
$ws.Range("J20").Value = 2826.4285
$ws.Range("K20").Value = 3614
$ws.Range("L20").Value = 2826.4285
$ws.Range("M20").Value = -3367
$ws.Range("N20").Value = -3320.4285

$ws.Range("H99").Value = 1099.6666
$ws.Range("I99").Value = 1099.6666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1099.6666
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 398.3334
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 1507
$ws.Range("I105").Value = 1507
$ws.Range("K105").Value = 1507
$ws.Range("M105").Value = 240

$ws.Range("H138").Value = 124500
$ws.Range("J138").Value = 124500
$ws.Range("L138").Value = 124500
$ws.Range("N138").Value = -134780

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2746.5293
$ws.Range("I31").Value = 1591.6923
$ws.Range("J31").Value = 6499.75
$ws.Range("K31").Value = 1591.6923
$ws.Range("L31").Value = 6499.75
$ws.Range("M31").Value = -1296.6923
$ws.Range("N31").Value = -7089.75

$ws.Range("H34").Value = 2746.5293
$ws.Range("I34").Value = 1591.6923
$ws.Range("J34").Value = 6499.75
$ws.Range("K34").Value = 1591.6923
$ws.Range("L34").Value = 6499.75
$ws.Range("M34").Value = -1389.6923
$ws.Range("N34").Value = -6903.75

$ws.Range("H43").Value = 11729.2
$ws.Range("J43").Value = 11729.2
$ws.Range("L43").Value = 11729.2
$ws.Range("N43").Value = -12097.2

$ws.Range("H96").Value = 32062
$ws.Range("J96").Value = 32062
$ws.Range("L96").Value = 32062
$ws.Range("N96").Value = -37554

$ws.Range("H97").Value = 33973.5
$ws.Range("J97").Value = 33973.5
$ws.Range("L97").Value = 33973.5
$ws.Range("N97").Value = -35955.5

$ws.Range("H101").Value = 11729.2
$ws.Range("J101").Value = 11729.2
$ws.Range("L101").Value = 11729.2
$ws.Range("N101").Value = -18219.2

$ws.Range("H105").Value = 2484.6
$ws.Range("I105").Value = 2484.3333
$ws.Range("J105").Value = 2485
$ws.Range("K105").Value = 2484.3333
$ws.Range("L105").Value = 2485
$ws.Range("M105").Value = -737.3332999999998
$ws.Range("N105").Value = -5979

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 777.6
$ws.Range("I12").Value = 679.3333
$ws.Range("J12").Value = 925
$ws.Range("K12").Value = 2037.9999
$ws.Range("L12").Value = 2775
$ws.Range("M12").Value = -1864.9999
$ws.Range("N12").Value = -3121

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H122").Value = 846.0833
$ws.Range("I122").Value = 516.25
$ws.Range("J122").Value = 1011
$ws.Range("K122").Value = 4646.25
$ws.Range("L122").Value = 9099
$ws.Range("M122").Value = -2196.25
$ws.Range("N122").Value = -13999

$ws.Range("H132").Value = 4888.3335
$ws.Range("I132").Value = 3333
$ws.Range("J132").Value = 6443.6665
$ws.Range("K132").Value = 29997
$ws.Range("L132").Value = 57992.9985
$ws.Range("M132").Value = -27467
$ws.Range("N132").Value = -63052.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3866.6667
$ws.Range("I102").Value = 3866.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3866.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2244.6667
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 1213.5
$ws.Range("J122").Value = 1213
$ws.Range("L122").Value = 3639
$ws.Range("N122").Value = -8539

$ws.Range("H126").Value = 9998.5
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 9998.5
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 29995.5
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -34935.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1705

$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1893

$ws.Range("H61").Value = 2266.3333
$ws.Range("I61").Value = 2266.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2266.3333
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2064.3333

$ws.Range("H113").Value = 2266.3333
$ws.Range("I113").Value = 2266.3333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2266.3333
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -96.33329999999978

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
